$d = $word.ActiveDocument

# Update the date title paragraph
$d.Content.Find.Execute("2023-08-18 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-19 Saturday", 2) | Out-Null

# Update each arithmetic expression cell in the table, addressed by position
# to safely handle the duplicate "26+57=" value that appears twice in the source.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "48+18="
$t.Cell(1, 2).Range.Text = "56-48="
$t.Cell(1, 3).Range.Text = "33+38="
$t.Cell(1, 4).Range.Text = "9+19="
$t.Cell(1, 5).Range.Text = "49+8="
$t.Cell(2, 1).Range.Text = "60-46="
$t.Cell(2, 2).Range.Text = "3+89="
$t.Cell(2, 3).Range.Text = "54-49="
$t.Cell(2, 4).Range.Text = "58+5="
$t.Cell(2, 5).Range.Text = "43-36="
$t.Cell(3, 1).Range.Text = "65+26="
$t.Cell(3, 2).Range.Text = "7+64="
$t.Cell(3, 3).Range.Text = "57+28="
$t.Cell(3, 4).Range.Text = "32-19="
$t.Cell(3, 5).Range.Text = "65-9="
$t.Cell(4, 1).Range.Text = "74-55="
$t.Cell(4, 2).Range.Text = "8+88="
$t.Cell(4, 3).Range.Text = "71-32="
$t.Cell(4, 4).Range.Text = "39+14="
$t.Cell(4, 5).Range.Text = "94-6="
$t.Cell(5, 1).Range.Text = "6+65="
$t.Cell(5, 2).Range.Text = "91-48="
$t.Cell(5, 3).Range.Text = "80-53="
$t.Cell(5, 4).Range.Text = "67+25="
$t.Cell(5, 5).Range.Text = "14+9="
$t.Cell(6, 1).Range.Text = "35+48="
$t.Cell(6, 2).Range.Text = "76+7="
$t.Cell(6, 3).Range.Text = "38+56="
$t.Cell(6, 4).Range.Text = "9+18="
$t.Cell(6, 5).Range.Text = "9+36="
$t.Cell(7, 1).Range.Text = "93-17="
$t.Cell(7, 2).Range.Text = "20-17="
$t.Cell(7, 3).Range.Text = "45+18="
$t.Cell(7, 4).Range.Text = "54+19="
$t.Cell(7, 5).Range.Text = "43-6="
$t.Cell(8, 1).Range.Text = "81-52="
$t.Cell(8, 2).Range.Text = "59+26="
$t.Cell(8, 3).Range.Text = "87-68="
$t.Cell(8, 4).Range.Text = "90-7="
$t.Cell(8, 5).Range.Text = "83-9="
$t.Cell(9, 1).Range.Text = "91-72="
$t.Cell(9, 2).Range.Text = "60-59="
$t.Cell(9, 3).Range.Text = "83-6="
$t.Cell(9, 4).Range.Text = "9+14="
$t.Cell(9, 5).Range.Text = "95-17="
$t.Cell(10, 1).Range.Text = "73-18="
$t.Cell(10, 2).Range.Text = "14+77="
$t.Cell(10, 3).Range.Text = "66-7="
$t.Cell(10, 4).Range.Text = "93-36="
$t.Cell(10, 5).Range.Text = "17+44="
$t.Cell(11, 1).Range.Text = "27+58="
$t.Cell(11, 2).Range.Text = "55+19="
$t.Cell(11, 3).Range.Text = "36+58="
$t.Cell(11, 4).Range.Text = "56+15="
$t.Cell(11, 5).Range.Text = "85-6="
$t.Cell(12, 1).Range.Text = "4+57="
$t.Cell(12, 2).Range.Text = "90-74="
$t.Cell(12, 3).Range.Text = "69+27="
$t.Cell(12, 4).Range.Text = "19+69="
$t.Cell(12, 5).Range.Text = "44+48="
$t.Cell(13, 1).Range.Text = "38+26="
$t.Cell(13, 2).Range.Text = "57+17="
$t.Cell(13, 3).Range.Text = "5+8="
$t.Cell(13, 4).Range.Text = "65+19="
$t.Cell(13, 5).Range.Text = "43-24="
$t.Cell(14, 1).Range.Text = "70-18="
$t.Cell(14, 2).Range.Text = "43-7="
$t.Cell(14, 3).Range.Text = "44+47="
$t.Cell(14, 4).Range.Text = "26+6="
$t.Cell(14, 5).Range.Text = "70-32="
$t.Cell(15, 1).Range.Text = "24-5="
$t.Cell(15, 2).Range.Text = "26-19="
$t.Cell(15, 3).Range.Text = "39+2="
$t.Cell(15, 4).Range.Text = "39+6="
$t.Cell(15, 5).Range.Text = "75-69="
$t.Cell(16, 1).Range.Text = "88-59="
$t.Cell(16, 2).Range.Text = "70-45="
$t.Cell(16, 3).Range.Text = "55-47="
$t.Cell(16, 4).Range.Text = "49+44="
$t.Cell(16, 5).Range.Text = "81-24="
$t.Cell(17, 1).Range.Text = "87+9="
$t.Cell(17, 2).Range.Text = "59+32="
$t.Cell(17, 3).Range.Text = "27+29="
$t.Cell(17, 4).Range.Text = "46+49="
$t.Cell(17, 5).Range.Text = "51-24="
$t.Cell(18, 1).Range.Text = "80-25="
$t.Cell(18, 2).Range.Text = "51-3="
$t.Cell(18, 3).Range.Text = "4+17="
$t.Cell(18, 4).Range.Text = "94-85="
$t.Cell(18, 5).Range.Text = "2+29="
$t.Cell(19, 1).Range.Text = "9+15="
$t.Cell(19, 2).Range.Text = "85-68="
$t.Cell(19, 3).Range.Text = "22-3="
$t.Cell(19, 4).Range.Text = "17+26="
$t.Cell(19, 5).Range.Text = "48+14="
$t.Cell(20, 1).Range.Text = "17+17="
$t.Cell(20, 2).Range.Text = "46+6="
$t.Cell(20, 3).Range.Text = "7+56="
$t.Cell(20, 4).Range.Text = "59+6="
$t.Cell(20, 5).Range.Text = "65-18="
